$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.054.35"
$ws.Range("E2").Value = "  -0.35%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.828.50"
$ws.Range("E3").Value = "  -0.29%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9989"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.65"
$ws.Range("E5").Value = "  -0.31%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6180"
$ws.Range("E6").Value = "  -6.93%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.43"
$ws.Range("E8").Value = "  +6.25%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07336"
$ws.Range("E9").Value = "  -1.16%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.2913"
$ws.Range("E10").Value = "  -0.82%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.67"
$ws.Range("E11").Value = "  +0.04%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07684"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.827.47"
$ws.Range("E13").Value = "  +0.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.963"
$ws.Range("E14").Value = "  -0.55%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6615"
$ws.Range("E15").Value = "  -1.21%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.81"
$ws.Range("E16").Value = "  -1.49%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008956"
$ws.Range("E17").Value = "  +6.60%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.019"
$ws.Range("E18").Value = "  -1.30%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.056.09"
$ws.Range("E19").Value = "  -0.32%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.075.58"
$ws.Range("E20").Value = "  +0.20%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "225.35"
$ws.Range("E21").Value = "  -0.85%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.35"
$ws.Range("E22").Value = "  -1.06%  "

$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.128"
$ws.Range("E24").Value = "  -0.80%  "

$ws.Range("E25").Value = "  +0.01%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.60"
$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.424"
$ws.Range("E27").Value = "  -2.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1352"
$ws.Range("E28").Value = "  -3.87%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "17.77"
$ws.Range("E29").Value = "  -1.01%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.495"
$ws.Range("E30").Value = "  -1.03%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.033"
$ws.Range("E31").Value = "  -0.37%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.047"
$ws.Range("E32").Value = "  -1.64%  "

$ws.Range("E33").Value = "  +0.37%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05249"
$ws.Range("E34").Value = "  -1.63%  "

$ws.Range("E35").Value = "  -1.87%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.149"
$ws.Range("E36").Value = "  +1.07%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.7304"
$ws.Range("E37").Value = "  -3.62%  "

$ws.Range("E38").Value = "  -0.92%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.295.21"
$ws.Range("E39").Value = "  +1.61%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.748"
$ws.Range("E40").Value = "  +0.47%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01782"
$ws.Range("E41").Value = "  -0.99%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.289"
$ws.Range("E42").Value = "  +5.25%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9000"
$ws.Range("E43").Value = "  -3.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9992"
$ws.Range("E44").Value = "  -0.58%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.84"
$ws.Range("E45").Value = "  -0.92%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.973.70"
$ws.Range("E46").Value = "  +0.03%  "

$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5114"
$ws.Range("E47").Value = "  -0.91%  "

$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "63.93"
$ws.Range("E48").Value = "  +0.92%  "

$ws.Range("E49").Value = "  -0.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.713"
$ws.Range("E50").Value = "  -3.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3961"
$ws.Range("E51").Value = "  -1.87%  "
